$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('D2').Value = '37.244.66'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '2.033.22'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '227.87'
$ws.Range('E5').Value = '  -2.66%  '
$ws.Range('D6').Value = '0.609'
$ws.Range('E6').Value = '  -2.60%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '55.65'
$ws.Range('E8').Value = '  -4.58%  '
$ws.Range('D9').Value = '0.383'
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('D10').Value = '0.0795'
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  -2.40%  '
$ws.Range('D12').Value = '2.338.48'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('D13').Value = '14.36'
$ws.Range('E13').Value = '  -5.58%  '
$ws.Range('D14').Value = '20.41'
$ws.Range('E14').Value = '  -3.82%  '
$ws.Range('D15').Value = '0.746'
$ws.Range('E15').Value = '  -4.09%  '
$ws.Range('D16').Value = '5.19'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').Value = '2.049.63'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '37.119.26'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = '6.01'
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').Value = '69.14'
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('D21').Value = '0.0₃0837'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '223.89'
$ws.Range('E22').Value = '  -2.58%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '2.38'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -5.31%  '
$ws.Range('D26').Value = '9.39'
$ws.Range('E26').Value = '  -3.55%  '
$ws.Range('D27').Value = '167.84'
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('E28').Value = '  -6.64%  '
$ws.Range('D29').Value = '18.78'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('D30').Value = '1.34'
$ws.Range('E30').Value = '  -3.68%  '
$ws.Range('E31').Value = '  -4.19%  '
$ws.Range('D32').Value = '4.49'
$ws.Range('E32').Value = '  -4.53%  '
$ws.Range('D33').Value = '0.0608'
$ws.Range('E33').Value = '  -4.42%  '
$ws.Range('D34').Value = '4.46'
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('E35').Value = '  -5.35%  '
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').Value = '3.16'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('D39').Value = '5.33'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('E40').Value = '  -7.57%  '
$ws.Range('D41').Value = '1.491.99'
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '2.85'
$ws.Range('E42').Value = '  -1.95%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').Value = '0.0935'
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('D44').Value = '95.19'
$ws.Range('E44').Value = '  -5.89%  '
$ws.Range('D45').Value = '16.54'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('D46').Value = '1.13'
$ws.Range('E46').Value = '  -6.17%  '
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  -5.15%  '
$ws.Range('D48').Value = '7.13'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('D50').Value = '3.69'
$ws.Range('E50').Value = '  -10.06%  '
$ws.Range('D51').Value = '2.226.63'
$ws.Range('E51').Value = '  -2.33%  '
